# Hotel Booking.pptx - "final mod on .ipynb"
# 1) Bump the cached date-field placeholders (Date Placeholder shapes on the
#    slide master / every slide layout, plus the Notes Master) by one day,
#    mirroring the real PowerPoint re-cache that happens when the deck is
#    reopened/saved on the next day.
# 2) Rewrite the body text + autofit of the "Results" slide (slide 4).

$p = $ppt.ActivePresentation

function Bump-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if (-not $sh.HasTextFrame) { continue }
        if ($sh.Name -notlike "Date Placeholder*") { continue }
        $tr = $sh.TextFrame.TextRange
        $txt = $tr.Text
        if ($txt -eq "2020-03-21") {
            $tr.Text = "2020-03-22"
        } elseif ($txt -eq "3/21/20") {
            $tr.Text = "3/22/20"
        }
    }
}

# Slide master
Bump-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout off the master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Bump-DatePlaceholder $layouts.Item($li).Shapes
}

# Notes master (datetimeFigureOut field)
Bump-DatePlaceholder $p.NotesMaster.Shapes

# --- Slide 4 ("Results") content rewrite ---
$s4 = $p.Slides.Item(4)
$content = $s4.Shapes.Item(2)
$tf = $content.TextFrame

$tf.TextRange.Text = "I used Binomial Logistic Regression to perform binomial classification.`r" + `
    "The main objective is to identify which feature(s) play important role in identifying weather a hotel booking will be cancelled or not. `r" + `
    "`r" + `
    "Based on 31 different features, there are several features that have significance on the final decision. For instance, people that has children and babies has tendencies to cancel their booking compared to those who doesn’t.`r" + `
    "`r" + `
    "For more detail, please look at the code and its output."

$tf.AutoSize = 2
